$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(11.359748798043237, 5.204022976389113, -11.419389052197639, 12.387197169402343)"
$ws.Range("C2").Value = "NCT(2.411391756868443, 1.687460250994135, 1.0838414493158748, 4.003545505552404)"
$ws.Range("D2").Value = "JSU(-0.846485049793761, 1.0937297549630995, 0.8609297517233687, 2.91409818972001)"
$ws.Range("E2").Value = "NCT(3.3850254793990135, 2.448048280435412, -3.3573525994004227, 4.623662666367011)"
